$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 4 (rows 4 & 5), pushing the existing
# rows 4..30 down to 6..32.
$ws.Rows.Item(4).EntireRow.Insert()
$ws.Rows.Item(4).EntireRow.Insert()

# Row 3: alias (A3) now reads "BRENDA ROMERO HERNANDEZ" while the
# nombre (B3) keeps "BRENDA GRISELDA ROMERO HERNANDEZ".
$ws.Range("A3").Value = "BRENDA ROMERO HERNANDEZ"
$ws.Range("B3").Value = "BRENDA GRISELDA ROMERO HERNANDEZ"

# New row 4: BRENDA ROMERO, DUI 054354354
$ws.Range("A4").Value = "BRENDA ROMERO"
$ws.Range("B4").Value = "BRENDA GRISELDA ROMERO HERNANDEZ"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "054354354"
$ws.Range("D4").Value = 122641732

# New row 5: GRISELDA HERNANDEZ, DUI 012345678 (with trailing space)
$ws.Range("A5").Value = "GRISELDA HERNANDEZ"
$ws.Range("B5").Value = "BRENDA GRISELDA ROMERO HERNANDEZ"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "012345678 "
$ws.Range("D5").Value = 122641733

# Match row 3's row height and the highlighted/bordered C-column format
# (yellow fill, grey border, wrap text) on the two new rows, without
# touching the values that were just written (PasteSpecial formats only).
$ws.Rows.Item(4).RowHeight = 13.8
$ws.Rows.Item(5).RowHeight = 13.8

$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection to B3, matching the saved cursor position.
$ws.Range("B3").Select()
